# Generate Report for Handoff
# Updates Priority (E) and Latest Handoff Datetime (H) for rows 4-7
# (the items that have now been handed off) on both the zh-cn and de-de
# worksheets, and refreshes the matching "Latest HO Xliff Generate Date"
# column on the Overview sheet.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# zh-cn sheet: rows 4-7 -> Priority becomes "ht", handoff datetime bumped
for ($r = 4; $r -le 7; $r++) {
    $zhcn.Cells.Item($r, 5).Value = "ht"
    $zhcn.Cells.Item($r, 8).Value = "2016-08-26 04:30:31"
}

# de-de sheet: rows 4-7 -> Priority becomes "ht", handoff datetime bumped
for ($r = 4; $r -le 7; $r++) {
    $dede.Cells.Item($r, 5).Value = "ht"
    $dede.Cells.Item($r, 8).Value = "2016-08-26 04:30:37"
}

# Overview sheet: rows 4-7 -> Latest HO Xliff Generate Date refreshed
# (this mirrors the new de-de handoff timestamp, which is the newer one)
for ($r = 4; $r -le 7; $r++) {
    $overview.Cells.Item($r, 7).Value = "2016-08-26 04:30:37"
}
